# Auto-generated: scheduled market-price data refresh for Atomos_Profits sheets
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) with
# freshly pulled market-board values for the affected Leve rows.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 6803032.5
$ws.Range("I19").Value = 11905033
$ws.Range("J19").Value = 364.44446
$ws.Range("K19").Value = 11905033
$ws.Range("L19").Value = 364.44446
$ws.Range("M19").Value = -11904858
$ws.Range("N19").Value = -714.4444599999999
$ws.Range("H116").Value = 6230
$ws.Range("I116").Value = 6640
$ws.Range("K116").Value = 6640
$ws.Range("M116").Value = -3198
$ws.Range("H123").Value = 36148.89
$ws.Range("J123").Value = 36148.89
$ws.Range("L123").Value = 36148.89
$ws.Range("N123").Value = -45948.89
$ws.Range("H134").Value = 24630
$ws.Range("J134").Value = 24630
$ws.Range("L134").Value = 24630
$ws.Range("N134").Value = -34770
$ws.Range("H138").Value = 2343.5483
$ws.Range("I138").Value = 1772.4286
$ws.Range("K138").Value = 5317.2858
$ws.Range("M138").Value = -177.2857999999997
$ws.Range("H141").Value = 423666.2
$ws.Range("I141").Value = 2520
$ws.Range("J141").Value = 564048.25
$ws.Range("K141").Value = 7560
$ws.Range("L141").Value = 1692144.75
$ws.Range("M141").Value = -2380
$ws.Range("N141").Value = -1702504.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6304.457
$ws.Range("I32").Value = 5307.3506
$ws.Range("J32").Value = 25498.75
$ws.Range("K32").Value = 5307.3506
$ws.Range("L32").Value = 25498.75
$ws.Range("M32").Value = -5020.3506
$ws.Range("N32").Value = -26072.75
$ws.Range("H97").Value = 897.0952
$ws.Range("I97").Value = 841.1875
$ws.Range("J97").Value = 1076
$ws.Range("K97").Value = 841.1875
$ws.Range("L97").Value = 1076
$ws.Range("M97").Value = -345.1875
$ws.Range("N97").Value = -2068

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 535.24
$ws.Range("I80").Value = 658
$ws.Range("J80").Value = 453.4
$ws.Range("K80").Value = 658
$ws.Range("L80").Value = 453.4
$ws.Range("M80").Value = 340
$ws.Range("N80").Value = -2449.4
$ws.Range("H83").Value = 535.24
$ws.Range("I83").Value = 658
$ws.Range("J83").Value = 453.4
$ws.Range("K83").Value = 3290
$ws.Range("L83").Value = 2267
$ws.Range("M83").Value = 1702
$ws.Range("N83").Value = -12251

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1861.5385
$ws.Range("I16").Value = 1100
$ws.Range("K16").Value = 1100
$ws.Range("M16").Value = -813
$ws.Range("H19").Value = 989.1667
$ws.Range("I19").Value = 187
$ws.Range("J19").Value = 5000
$ws.Range("K19").Value = 187
$ws.Range("L19").Value = 5000
$ws.Range("M19").Value = -17
$ws.Range("N19").Value = -5340
$ws.Range("H21").Value = 58735.75
$ws.Range("I21").Value = 37456.5
$ws.Range("J21").Value = 80015
$ws.Range("K21").Value = 37456.5
$ws.Range("L21").Value = 80015
$ws.Range("M21").Value = -37221.5
$ws.Range("N21").Value = -80485
$ws.Range("H24").Value = 989.1667
$ws.Range("I24").Value = 187
$ws.Range("J24").Value = 5000
$ws.Range("K24").Value = 187
$ws.Range("L24").Value = 5000
$ws.Range("M24").Value = -17
$ws.Range("N24").Value = -5340
$ws.Range("H105").Value = 3175
$ws.Range("I105").Value = 2444.4443
$ws.Range("K105").Value = 2444.4443
$ws.Range("M105").Value = -697.4443000000001
$ws.Range("H113").Value = 1861.5385
$ws.Range("I113").Value = 1100
$ws.Range("K113").Value = 1100
$ws.Range("M113").Value = 1070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1834.2142
$ws.Range("I97").Value = 999
$ws.Range("J97").Value = 1898.4615
$ws.Range("K97").Value = 2997
$ws.Range("L97").Value = 5695.3845
$ws.Range("M97").Value = -2501
$ws.Range("N97").Value = -6687.3845
$ws.Range("H126").Value = 2666.6667
$ws.Range("J126").Value = 2666.6667
$ws.Range("L126").Value = 8000.000100000001
$ws.Range("N126").Value = -17880.0001
$ws.Range("H132").Value = 2742.5715
$ws.Range("I132").Value = 1349.75
$ws.Range("J132").Value = 4599.6665
$ws.Range("K132").Value = 12147.75
$ws.Range("L132").Value = 41396.9985
$ws.Range("M132").Value = -9617.75
$ws.Range("N132").Value = -46456.9985

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 9816
$ws.Range("J5").Value = 13850
$ws.Range("L5").Value = 13850
$ws.Range("N5").Value = -14074
$ws.Range("H24").Value = 20787.715
$ws.Range("J24").Value = 20787.715
$ws.Range("L24").Value = 20787.715
$ws.Range("N24").Value = -21133.715
$ws.Range("H102").Value = 36419.266
$ws.Range("J102").Value = 95171.73
$ws.Range("L102").Value = 95171.73
$ws.Range("N102").Value = -98415.73

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3227099.2
$ws.Range("I7").Value = 5000771
$ws.Range("K7").Value = 5000771
$ws.Range("M7").Value = -5000659
$ws.Range("H14").Value = 139003.89
$ws.Range("I14").Value = 836734.7
$ws.Range("J14").Value = 8179.375
$ws.Range("K14").Value = 836734.7
$ws.Range("L14").Value = 8179.375
$ws.Range("M14").Value = -836562.7
$ws.Range("N14").Value = -8523.375
$ws.Range("H93").Value = 1608.0834
$ws.Range("I93").Value = 810.3333
$ws.Range("K93").Value = 810.3333
$ws.Range("M93").Value = 437.6667
$ws.Range("H126").Value = 3227099.2
$ws.Range("I126").Value = 5000771
$ws.Range("K126").Value = 15002313
$ws.Range("M126").Value = -14999843

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 26603.6
$ws.Range("I24").Value = 10009
$ws.Range("J24").Value = 30752.25
$ws.Range("K24").Value = 10009
$ws.Range("L24").Value = 30752.25
$ws.Range("M24").Value = -9779
$ws.Range("N24").Value = -31212.25
$ws.Range("H81").Value = 1576.3043
$ws.Range("I81").Value = 1299
$ws.Range("K81").Value = 2598
$ws.Range("M81").Value = -1537
$ws.Range("H84").Value = 1576.3043
$ws.Range("I84").Value = 1299
$ws.Range("K84").Value = 12990
$ws.Range("M84").Value = -7686
$ws.Range("H109").Value = 37233.332
$ws.Range("J109").Value = 37233.332
$ws.Range("L109").Value = 37233.332
$ws.Range("N109").Value = -40007.332
$ws.Range("H126").Value = 4002252.8
$ws.Range("I126").Value = 1425.6875
$ws.Range("J126").Value = 1425.6875
$ws.Range("K126").Value = 4277.0625
$ws.Range("M126").Value = -1807.0625
